$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Add "(Graphical user interface)" paragraph right after the "Usability"
#    heading paragraph in the left-hand header cell.
# ---------------------------------------------------------------------------
$usabilityPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Usability`r`a") {
        $usabilityPara = $p
        break
    }
}
if ($usabilityPara -ne $null) {
    $usabilityPara.Range.InsertParagraphAfter()
    $newPara = $usabilityPara.Next()
    $newRng = $newPara.Range
    $newRng.Text = "(Graphical user interface)"
    $newRng.Font.Name = "Times New Roman"
    $newRng.Font.Bold = $true
    $newRng.Font.Size = 14
    $newRng.Font.SizeBi = 14
}

Write-Host "step1 done"

# ---------------------------------------------------------------------------
# 2. Expand the "HTML and css ... Db system :" bullet: insert a long new
#    sentence between "GUI systems" and "Db system :".
# ---------------------------------------------------------------------------
$newSentence = " are friendly to user experience. Categorized items are easier to find what customers need. Items are always posted in order of popularity from the top on the list of items. Seller" + [char]0x2019 + "s hash tag also gives more opportunities to sell his items. Db system :"
$found2 = $d.Content.Find.Execute(" Db system :", $true, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2)
Write-Host "step2 found=$found2"

# ---------------------------------------------------------------------------
# 3. Collapse the split "Access al/l /user/ who buy" runs into one run of
#    text "Access all user who buy".
# ---------------------------------------------------------------------------
$found3 = $d.Content.Find.Execute("Access al" + "l " + "user" + " who buy", $true, $false, $false, $false, $false, $true, 1, $false, "Access all user who buy", 2)
Write-Host "step3 found=$found3"

# ---------------------------------------------------------------------------
# 4. Append " Login id is based on Email." after the "Account function..."
#    sentence.
# ---------------------------------------------------------------------------
$acctText = "Account function: create, remove and edit(allow changing pw and card account  payment method)."
$found4 = $d.Content.Find.Execute($acctText, $true, $false, $false, $false, $false, $true, 1, $false, $acctText + " Login id is based on Email.", 2)
Write-Host "step4 found=$found4"

# ---------------------------------------------------------------------------
# 5. "Password and cvv code..." -> "Password and credit card cvv code..."
# ---------------------------------------------------------------------------
$found5 = $d.Content.Find.Execute("Password and cvv code", $true, $false, $false, $false, $false, $true, 1, $false, "Password and credit card cvv code", 2)
Write-Host "step5 found=$found5"
